$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.420.59'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.62%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.642.01'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.42%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.07'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.535'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.07%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.04'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0608'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.11%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0892'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.875.61'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.38%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.645.80'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.62%  '
$ws.Range('E14').Value = '  -2.98%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.24'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.90%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.399.01'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '228.18'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -7.62%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0718'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.46'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.07%  '
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('E22').Value = '  -3.80%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('E24').Value = '  -0.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '148.02'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.32%  '
$ws.Range('E26').Value = '  +1.52%  '
$ws.Range('E27').Value = '  -3.48%  '
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.53'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.22%  '
$ws.Range('E30').Value = '  -5.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0485'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.05%  '
$ws.Range('E32').Value = '  -2.12%  '
$ws.Range('E33').Value = '  -0.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.409.37'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.74%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('E36').Value = '  -0.28%  '
$ws.Range('E37').Value = '  -2.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.880'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.11%  '
$ws.Range('E39').Value = '  -3.27%  '
$ws.Range('E40').Value = '  +1.22%  '
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('E42').Value = '  -1.54%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.46'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.96%  '
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '64.60'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.23%  '
$ws.Range('E46').Value = '  +0.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.786.04'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.23%  '
$ws.Range('E48').Value = '  -3.90%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '87.34'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.14%  '
$ws.Range('E50').Value = '  -3.30%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0986'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.61%  '
